$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 4355
$ws.Range("I2").Value = 2616.6667
$ws.Range("J2").Value = 5224.1665
$ws.Range("K2").Value = 2616.6667
$ws.Range("L2").Value = 5224.1665
$ws.Range("M2").Value = -2503.6667
$ws.Range("N2").Value = -5450.1665
$ws.Range("H40").Value = 2578.4
$ws.Range("I40").Value = 3397.5
$ws.Range("J40").Value = 2032.3334
$ws.Range("K40").Value = 3397.5
$ws.Range("L40").Value = 2032.3334
$ws.Range("M40").Value = -3222.5
$ws.Range("N40").Value = -2382.3334
$ws.Range("H55").Value = 2498.8
$ws.Range("I55").Value = 247.5
$ws.Range("J55").Value = 3999.6667
$ws.Range("K55").Value = 247.5
$ws.Range("L55").Value = 3999.6667
$ws.Range("M55").Value = -33.5
$ws.Range("N55").Value = -4427.6667
$ws.Range("H98").Value = 3432.2
$ws.Range("I98").Value = 3358.7273
$ws.Range("J98").Value = 3634.25
$ws.Range("K98").Value = 3358.7273
$ws.Range("L98").Value = 3634.25
$ws.Range("M98").Value = -1860.7273
$ws.Range("N98").Value = -6630.25
$ws.Range("H107").Value = 19609000
$ws.Range("I107").Value = 1229
$ws.Range("J107").Value = 333333340
$ws.Range("K107").Value = 1229
$ws.Range("L107").Value = 333333340
$ws.Range("M107").Value = 691
$ws.Range("N107").Value = -333337180
$ws.Range("H108").Value = 128592.336
$ws.Range("J108").Value = 128592.336
$ws.Range("L108").Value = 128592.336
$ws.Range("N108").Value = -136272.336
$ws.Range("H109").Value = 97543.8
$ws.Range("I109").Value = 40819
$ws.Range("J109").Value = 111725
$ws.Range("K109").Value = 40819
$ws.Range("L109").Value = 111725
$ws.Range("M109").Value = -39432
$ws.Range("N109").Value = -114499
$ws.Range("H113").Value = 3879.6
$ws.Range("J113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("N113").ClearContents()
$ws.Range("H116").Value = 5484.75
$ws.Range("I116").Value = 5429.6665
$ws.Range("J116").Value = 5517.8
$ws.Range("K116").Value = 5429.6665
$ws.Range("L116").Value = 5517.8
$ws.Range("M116").Value = -1987.6665
$ws.Range("N116").Value = -12401.8
$ws.Range("H122").Value = 3432.2
$ws.Range("I122").Value = 3358.7273
$ws.Range("J122").Value = 3634.25
$ws.Range("K122").Value = 10076.1819
$ws.Range("L122").Value = 10902.75
$ws.Range("M122").Value = -7626.1819
$ws.Range("N122").Value = -15802.75
$ws.Range("H126").Value = 130000
$ws.Range("J126").Value = 130000
$ws.Range("L126").Value = 130000
$ws.Range("N126").Value = -139880
$ws.Range("H127").Value = 9210.666999999999
$ws.Range("I127").Value = 1316
$ws.Range("J127").Value = 25000
$ws.Range("K127").Value = 3948
$ws.Range("L127").Value = 75000
$ws.Range("M127").Value = 1012
$ws.Range("N127").Value = -84920
$ws.Range("H135").Value = 981.6667
$ws.Range("I135").Value = 975
$ws.Range("K135").Value = 8775
$ws.Range("M135").Value = -6240

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H34").Value = 254750
$ws.Range("I34").Value = 173000
$ws.Range("K34").Value = 173000
$ws.Range("M34").Value = -172729
$ws.Range("H45").Value = 144949
$ws.Range("I45").Value = 201928.6
$ws.Range("K45").Value = 201928.6
$ws.Range("M45").Value = -201551.6
$ws.Range("H110").Value = 970.069
$ws.Range("I110").Value = 894.6
$ws.Range("J110").Value = 1441.75
$ws.Range("K110").Value = 894.6
$ws.Range("L110").Value = 1441.75
$ws.Range("M110").Value = 1150.4
$ws.Range("N110").Value = -5531.75
$ws.Range("H132").Value = 1750.1666
$ws.Range("I132").Value = 1231.5
$ws.Range("J132").Value = 3306.1667
$ws.Range("K132").Value = 3694.5
$ws.Range("L132").Value = 9918.500100000001
$ws.Range("M132").Value = -1164.5
$ws.Range("N132").Value = -14978.5001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value = 35000
$ws.Range("J81").Value = 35000
$ws.Range("L81").Value = 35000
$ws.Range("N81").Value = -37122
$ws.Range("H84").Value = 35000
$ws.Range("J84").Value = 35000
$ws.Range("L84").Value = 105000
$ws.Range("N84").Value = -115608
$ws.Range("H110").Value = 100333.336
$ws.Range("J110").Value = 100333.336
$ws.Range("L110").Value = 100333.336
$ws.Range("N110").Value = -108513.336
$ws.Range("H135").Value = 78569.57000000001
$ws.Range("J135").Value = 78569.57000000001
$ws.Range("L135").Value = 78569.57000000001
$ws.Range("N135").Value = -88709.57000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H55").Value = 15050
$ws.Range("I55").Value = 100
$ws.Range("J55").Value = 30000
$ws.Range("K55").Value = 100
$ws.Range("L55").Value = 30000
$ws.Range("M55").Value = 215
$ws.Range("N55").Value = -30630
$ws.Range("H58").Value = 2028.6
$ws.Range("I58").Value = 1916.1666
$ws.Range("J58").Value = 2103.5557
$ws.Range("K58").Value = 1916.1666
$ws.Range("L58").Value = 2103.5557
$ws.Range("M58").Value = -1713.1666
$ws.Range("N58").Value = -2509.5557
$ws.Range("H107").Value = 2249.7
$ws.Range("I107").Value = 2342.8572
$ws.Range("J107").Value = 2032.3334
$ws.Range("K107").Value = 2342.8572
$ws.Range("L107").Value = 2032.3334
$ws.Range("M107").Value = -422.8571999999999
$ws.Range("N107").Value = -5872.3334
$ws.Range("H132").Value = 2391.647
$ws.Range("I132").Value = 1986.4
$ws.Range("K132").Value = 5959.200000000001
$ws.Range("M132").Value = -3429.200000000001
$ws.Range("H134").Value = 2215.3333
$ws.Range("I134").Value = 1911.5714
$ws.Range("J134").Value = 2822.8572
$ws.Range("K134").Value = 5734.7142
$ws.Range("L134").Value = 8468.571599999999
$ws.Range("M134").Value = -3199.7142
$ws.Range("N134").Value = -13538.5716
$ws.Range("H136").Value = 2028.6
$ws.Range("I136").Value = 1916.1666
$ws.Range("J136").Value = 2103.5557
$ws.Range("K136").Value = 5748.4998
$ws.Range("L136").Value = 6310.6671
$ws.Range("M136").Value = -3198.4998
$ws.Range("N136").Value = -11410.6671

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 2700.1667
$ws.Range("I68").Value = 2238
$ws.Range("J68").Value = 3162.3333
$ws.Range("K68").Value = 6714
$ws.Range("L68").Value = 9486.999899999999
$ws.Range("M68").Value = -5903
$ws.Range("N68").Value = -11108.9999
$ws.Range("H71").Value = 2700.1667
$ws.Range("I71").Value = 2238
$ws.Range("J71").Value = 3162.3333
$ws.Range("K71").Value = 20142
$ws.Range("L71").Value = 28460.9997
$ws.Range("M71").Value = -16086
$ws.Range("N71").Value = -36572.9997
$ws.Range("H80").Value = 3299
$ws.Range("I80").Value = 2098
$ws.Range("J80").Value = 4099.6665
$ws.Range("K80").Value = 6294
$ws.Range("L80").Value = 12298.9995
$ws.Range("M80").Value = -5358
$ws.Range("N80").Value = -14170.9995
$ws.Range("H83").Value = 3299
$ws.Range("I83").Value = 2098
$ws.Range("J83").Value = 4099.6665
$ws.Range("K83").Value = 18882
$ws.Range("L83").Value = 36896.9985
$ws.Range("M83").Value = -14202
$ws.Range("N83").Value = -46256.9985
$ws.Range("H92").Value = 600
$ws.Range("I92").Value = 600
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 1800
$ws.Range("L92").Value = 0
$ws.Range("M92").Value = -552
$ws.Range("N92").ClearContents()
$ws.Range("H94").Value = 6620.6665
$ws.Range("J94").Value = 6742.4
$ws.Range("L94").Value = 20227.2
$ws.Range("N94").Value = -21579.2
$ws.Range("H129").Value = 5259.533
$ws.Range("I129").Value = 1310.4286
$ws.Range("K129").Value = 3931.2858
$ws.Range("M129").Value = 1068.7142

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H44").Value = 24248
$ws.Range("J44").Value = 29664
$ws.Range("L44").Value = 29664
$ws.Range("N44").Value = -30856
$ws.Range("H97").Value = 143508.08
$ws.Range("I97").Value = 100300.3
$ws.Range("J97").Value = 251527.5
$ws.Range("K97").Value = 100300.3
$ws.Range("L97").Value = 251527.5
$ws.Range("M97").Value = -99804.3
$ws.Range("N97").Value = -252519.5
$ws.Range("H102").Value = 17242902
$ws.Range("I102").Value = 20834758
$ws.Range("K102").Value = 20834758
$ws.Range("M102").Value = -20833136
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents()
$ws.Range("H111").Value = 134000
$ws.Range("J111").Value = 134000
$ws.Range("L111").Value = 134000
$ws.Range("N111").Value = -140134
$ws.Range("H132").Value = 697477.7
$ws.Range("I132").Value = 6012.3477
$ws.Range("K132").Value = 18037.0431
$ws.Range("M132").Value = -15507.0431

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 8284.809999999999
$ws.Range("J46").Value = 3665.7222
$ws.Range("L46").Value = 3665.7222
$ws.Range("N46").Value = -4041.7222
$ws.Range("H82").Value = 555.3333
$ws.Range("J82").Value = 462.5
$ws.Range("L82").Value = 462.5
$ws.Range("N82").Value = -1184.5
$ws.Range("H85").Value = 555.3333
$ws.Range("J85").Value = 462.5
$ws.Range("L85").Value = 462.5
$ws.Range("N85").Value = -2958.5
$ws.Range("H93").Value = 2080.2856
$ws.Range("I93").Value = 2080.2856
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 2080.2856
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = -832.2856000000002
$ws.Range("N93").ClearContents()
$ws.Range("H102").Value = 75000
$ws.Range("J102").Value = 75000
$ws.Range("L102").Value = 75000
$ws.Range("N102").Value = -81490
$ws.Range("H132").Value = 3289.1155
$ws.Range("I132").Value = 2616.2273
$ws.Range("K132").Value = 7848.6819
$ws.Range("M132").Value = -5318.6819

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 3391.2
$ws.Range("I126").Value = 3310.16
$ws.Range("K126").Value = 9930.48
$ws.Range("M126").Value = -7460.48
$ws.Range("H132").Value = 2832.1052
$ws.Range("J132").Value = 3362.3845
$ws.Range("L132").Value = 10087.1535
$ws.Range("N132").Value = -15147.1535

Write-Output "Applied all market-data cell updates."